$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1, using the same formatting as the existing headers (A1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Fill boolean FALSE values for rows 2-12, columns F:H
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 6).Value = $false
    $ws.Cells.Item($row, 7).Value = $false
    $ws.Cells.Item($row, 8).Value = $false
}
